# Select cell C7 and clear its value (keeping cell formatting/style intact),
# mirroring the user selecting C7 and pressing Delete.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("C7")
$cell.Select()
$cell.ClearContents()
